$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the dead "Homework help" (lab04, row 25) and "Pandas" (lab05, row 26)
# links/topics from column D of the syllabus.
$ws.Range("D25").ClearContents()
$ws.Range("D26").ClearContents()

# Update the active selection to match the edited cell.
$null = $ws.Range("D26").Select()
